# Apply the "adding averages and more checks" edit to the workbook.
$wb = $excel.ActiveWorkbook
$wsTrain = $wb.Worksheets.Item(1)   # "Training Dashboard"
$wsExam  = $wb.Worksheets.Item(2)   # "Exam Dashboard"

# -----------------------------------------------------------------
# Helper: write a date-like literal string into a cell without Excel
# re-interpreting it as a date serial number. We do this by writing a
# formula that evaluates to the literal text, then freezing the
# result back into a plain value via Copy / PasteSpecial (values
# only). This preserves the cell's existing number format / style.
# (Positional parameters only -- named binding isn't reliable here.)
# -----------------------------------------------------------------
function Set-LiteralText($range, $text) {
    $range.Formula = '="' + $text + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# -----------------------------------------------------------------
# Training Dashboard sheet: PERIOD TO EXPIRE (col H) drops by 8 days
# and LAST UPDATE (col I) moves from 08-Sep-2025 to 16-Sep-2025 for
# every data row (3-11).
# -----------------------------------------------------------------
$rows = @(
    @{ Row = 3;  H = 482 },
    @{ Row = 4;  H = 364 },
    @{ Row = 5;  H = 448 },
    @{ Row = 6;  H = 588 },
    @{ Row = 7;  H = 423 },
    @{ Row = 8;  H = 588 },
    @{ Row = 9;  H = 86  },
    @{ Row = 10; H = 237 },
    @{ Row = 11; H = 331 }
)

foreach ($r in $rows) {
    $wsTrain.Cells.Item($r.Row, 8).Value = $r.H
    Set-LiteralText $wsTrain.Cells.Item($r.Row, 9) "16-Sep-2025"
}

$wsTrain.Range("I1").Select() | Out-Null
$excel.CutCopyMode = $false

# -----------------------------------------------------------------
# Exam Dashboard sheet: widen the COMMENTS column and update its
# remark text from "OK" to "date is valid" for the first three rows.
# -----------------------------------------------------------------
$wsExam.Range("E1").ColumnWidth = 14.1666666667

$wsExam.Range("E3").Value = "date is valid"
$wsExam.Range("E4").Value = "date is valid"
$wsExam.Range("E5").Value = "date is valid"

# -----------------------------------------------------------------
# Header-row styling: the bold header font gains a white color (it
# already sits on the dark blue fill), and the big bold title font
# collapses onto that same bold/white font (losing its 14pt size).
# -----------------------------------------------------------------
$wsTrain.Range("A2:K2").Font.Color = 16777215
$wsExam.Range("A2:G2").Font.Color = 16777215

$wsTrain.Range("A1").Font.Color = 16777215
$wsTrain.Range("A1").Font.Size = 11

$wsExam.Range("A1").Font.Color = 16777215
$wsExam.Range("A1").Font.Size = 11
